# Automatic tracker update: fill in "resultado" (G) and "profit" (H)
# columns for rows whose bet outcome has just become known.
#
# Business rule (inferred from existing populated rows):
#   - "Acierto" (hit)  -> profit = cuota (F) - 1
#   - "Fallo"   (miss) -> profit = -1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 120; Resultado = "Acierto"; Profit = 0.73 },
    @{ Row = 126; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 127; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 128; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 140; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 142; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 148; Resultado = "Fallo";   Profit = -1   },
    @{ Row = 149; Resultado = "Acierto"; Profit = 0.83 },
    @{ Row = 150; Resultado = "Acierto"; Profit = 0.36 },
    @{ Row = 155; Resultado = "Fallo";   Profit = -1   }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 7).Value = $u.Resultado   # column G: resultado
    $ws.Cells.Item($r, 8).Value = $u.Profit      # column H: profit
}
